$wb = $excel.ActiveWorkbook

# --- Update the "Metrics" sheet values (B2:B13) ---
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 176018.44
$metrics.Range("B3").Value = 151133.37000000002
$metrics.Range("B4").Value = 54015.549999999996
$metrics.Range("B5").Value = 7201
$metrics.Range("B6").Value = 5378725.5500000007
$metrics.Range("B7").Value = 4551486.330000001
$metrics.Range("B8").Value = 1585972.4300000004
$metrics.Range("B9").Value = 209908
$metrics.Range("B10").Value = 33844106.539999992
$metrics.Range("B11").Value = 31826761.490000002
$metrics.Range("B12").Value = 11867694.469999995
$metrics.Range("B13").Value = 1307538

# Update selection on Metrics sheet to D14
$metrics.Range("D14").Select()

# --- Update the "today" sheet selection ---
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("D7").Select()

$excel.Calculate()
